$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1 (03:05 -> 04:05)
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 04:05"

# Row 51: Corea del Sur
$ws.Range("A51").Value = "Corea del Sur"
$ws.Range("B51").Value = 11503
$ws.Range("C51").Value = 35
$ws.Range("D51").Value = 10422
$ws.Range("E51").Value = 810
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 271

# Row 66: Australia
$ws.Range("A66").Value = "Australia"
$ws.Range("B66").Value = 7202
$ws.Range("C66").Value = 7
$ws.Range("D66").Value = 6618
$ws.Range("E66").Value = 481
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 103

# Row 72: Guatemala
$ws.Range("A72").Value = "Guatemala"
$ws.Range("B72").Value = 5087
$ws.Range("C72").Value = 348
$ws.Range("D72").Value = 735
$ws.Range("E72").Value = 4244
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 6
$ws.Range("H72").Value = 108

# Row 73: Sudan
$ws.Range("A73").Value = "Sudan"
$ws.Range("B73").Value = 5026
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 1423
$ws.Range("E73").Value = 3317
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 286

# Row 91: Haiti
$ws.Range("A91").Value = "Haiti"
$ws.Range("B91").Value = 2124
$ws.Range("C91").Value = 259
$ws.Range("D91").Value = 24
$ws.Range("E91").Value = 2056
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 44

# Row 92: Cuba
$ws.Range("A92").Value = "Cuba"
$ws.Range("B92").Value = 2045
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 1809
$ws.Range("E92").Value = 153
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 83

# Row 93: Somalia
$ws.Range("A93").Value = "Somalia"
$ws.Range("B93").Value = 1976
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 348
$ws.Range("E93").Value = 1550
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 78

# Row 94: Kenia
$ws.Range("A94").Value = "Kenia"
$ws.Range("B94").Value = 1962
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 478
$ws.Range("E94").Value = 1420
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 64

# Row 95: Estonia
$ws.Range("A95").Value = "Estonia"
$ws.Range("B95").Value = 1869
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 1624
$ws.Range("E95").Value = 177
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 68

# Row 201: Santa Lucia
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("B201").Value = 18
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 18
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

# Row 202: Belice
$ws.Range("A202").Value = "Belice"
$ws.Range("B202").Value = 18
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 16
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 2

# Row 210: Montserrat
$ws.Range("A210").Value = "Montserrat"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 10
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 1

# Row 211: Seychelles
$ws.Range("A211").Value = "Seychelles"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
